$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) sometimes holds values that look like plain
# numbers (e.g. "605.65"). The source data stores these as literal text
# (European/grouped formatting such as "66.205.57" elsewhere in the same
# column), so force Text format before writing any value that Excel would
# otherwise auto-convert to a number, to avoid float rounding artifacts
# and keep the cell type consistent with the rest of the column.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
}

# --- Row swap/update for Maker / TheGraph / Kaspa (rows 42-44) ---
$ws.Range("B42").Value = "TheGraph"
$ws.Range("C42").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
Set-TextValue $ws.Range("D42") "0.297"
$ws.Range("E42").Value = "  +2.87%  "

$ws.Range("B43").Value = "Kaspa"
$ws.Range("C43").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue $ws.Range("D43") "0.120"
$ws.Range("E43").Value = "  -0.26%  "

$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue $ws.Range("D44") "3.018.99"
$ws.Range("E44").Value = "  -1.76%  "

# --- Price / Volume(1h) updates for remaining rows ---
$ws.Range("D2").Value = "66.205.57"
$ws.Range("E2").Value = "  +1.65%  "
$ws.Range("D3").Value = "3.224.79"
$ws.Range("E3").Value = "  +1.33%  "
$ws.Range("E4").Value = "  -0.07%  "
Set-TextValue $ws.Range("D5") "605.65"
$ws.Range("E5").Value = "  +4.57%  "
Set-TextValue $ws.Range("D6") "155.26"
$ws.Range("E6").Value = "  +2.72%  "
Set-TextValue $ws.Range("D7") "1.00"
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "3.222.83"
$ws.Range("E8").Value = "  +1.29%  "
$ws.Range("E9").Value = "  +0.74%  "
Set-TextValue $ws.Range("D10") "0.162"
$ws.Range("E10").Value = "  -0.99%  "
$ws.Range("E11").Value = "  -0.63%  "
$ws.Range("E12").Value = "  +0.91%  "
$ws.Range("E13").Value = "  -0.53%  "
Set-TextValue $ws.Range("D14") "38.99"
$ws.Range("E14").Value = "  +2.22%  "
$ws.Range("D15").Value = "3.744.72"
$ws.Range("E15").Value = "  +0.97%  "
Set-TextValue $ws.Range("D16") "7.51"
$ws.Range("E16").Value = "  +4.40%  "
$ws.Range("D17").Value = "66.337.35"
$ws.Range("E17").Value = "  +1.70%  "
$ws.Range("D18").Value = "3.234.53"
$ws.Range("E18").Value = "  +1.36%  "
$ws.Range("E19").Value = "  +0.68%  "
Set-TextValue $ws.Range("D20") "514.48"
$ws.Range("E20").Value = "  -0.10%  "
Set-TextValue $ws.Range("D21") "15.86"
$ws.Range("E21").Value = "  +6.21%  "
Set-TextValue $ws.Range("D22") "0.740"
$ws.Range("E22").Value = "  +0.60%  "
Set-TextValue $ws.Range("D23") "15.32"
$ws.Range("E23").Value = "  -0.68%  "
$ws.Range("E24").Value = "  +2.23%  "
Set-TextValue $ws.Range("D25") "85.70"
$ws.Range("E25").Value = "  +0.35%  "
Set-TextValue $ws.Range("D26") "0.999"
$ws.Range("E26").Value = "  -0.11%  "
Set-TextValue $ws.Range("D27") "3.05"
$ws.Range("E27").Value = "  +4.07%  "
Set-TextValue $ws.Range("D28") "9.29"
$ws.Range("E28").Value = "  +2.73%  "
$ws.Range("E29").Value = "  +2.66%  "
Set-TextValue $ws.Range("D30") "2.89"
$ws.Range("E30").Value = "  +3.88%  "
Set-TextValue $ws.Range("D31") "6.90"
$ws.Range("E31").Value = "  +10.34%  "
Set-TextValue $ws.Range("D32") "28.32"
$ws.Range("E32").Value = "  +0.72%  "
$ws.Range("E33").Value = "  +1.09%  "
$ws.Range("E34").Value = "  +0.06%  "
$ws.Range("E35").Value = "  +0.86%  "
$ws.Range("E36").Value = "  -0.13%  "
Set-TextValue $ws.Range("D37") "0.0924"
$ws.Range("E37").Value = "  +2.36%  "
Set-TextValue $ws.Range("D38") "492.75"
$ws.Range("E38").Value = "  +2.73%  "
$ws.Range("E39").Value = "  +0.05%  "
Set-TextValue $ws.Range("D40") "3.04"
$ws.Range("E40").Value = "  -3.35%  "
Set-TextValue $ws.Range("D41") "8.89"
$ws.Range("E41").Value = "  +2.38%  "
Set-TextValue $ws.Range("D45") "2.52"
$ws.Range("E45").Value = "  +4.03%  "
$ws.Range("D46").Value = "0.0₃0651"
$ws.Range("E46").Value = "  +7.12%  "
Set-TextValue $ws.Range("D47") "29.29"
$ws.Range("E47").Value = "  -0.10%  "
$ws.Range("E48").Value = "  +0.08%  "
$ws.Range("E49").Value = "  +0.21%  "
Set-TextValue $ws.Range("D50") "2.35"
$ws.Range("E50").Value = "  +2.95%  "
Set-TextValue $ws.Range("D51") "120.04"
$ws.Range("E51").Value = "  -0.55%  "
